$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the quote-prefixed text formatting used by row 5 down onto the two
# new rows before filling them in
$ws.Range("B5:C5").Copy() | Out-Null
$ws.Range("B6:C7").PasteSpecial(-4122) | Out-Null

# Row 5 (PREPROD): replace the old smart-folder values with the new ones
$ws.Range("C5").Value = "'2026909"
$ws.Range("B5").Value = "'0420172008629 "

# Row 6 (new): PREPROD smart folder to perform the cancellation itself
$ws.Range("A6").Value = "PREPROD"
$ws.Range("B6").Value = "'1120170200942 "
$ws.Range("C6").Value = "'0200212"

# Row 7 (new): PREPROD smart folder to obtain the cancellation number
$ws.Range("A7").Value = "PREPROD"
$ws.Range("B7").Value = "'1220170301442 "
$ws.Range("C7").Value = "'0300190 "

$ws.Range("C8").Select() | Out-Null
